$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already has one "table" (rows 3-5: spacer row, header row, data
# row) describing the AlexNet run. This edit duplicates that block into
# rows 6-8 for a second (ResNet50) run, leaving rows 9 onward as they were
# (row 9's A:I cells - which used to belong to nothing in particular - are
# removed since the duplicated table now ends at row 8).

# Row 6: blank "spacer" row - copy formatting only from row 3's pattern.
$ws.Range("A3:I3").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)

# Row 7: header labels (Network, Solver, Max epochs, ...) - copy values then
# formatting from row 4.
$ws.Range("A4:I4").Copy()
$ws.Range("A7:I7").PasteSpecial()
$ws.Range("A4:I4").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)

# Row 8: data row - copy values then formatting from row 5 (the AlexNet
# result row) as a starting point.
$ws.Range("A5:I5").Copy()
$ws.Range("A8:I8").PasteSpecial()
$ws.Range("A5:I5").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)

# Row 8 overrides: this run is for ResNet50, and the attempt didn't actually
# finish successfully, so the measured results (Accuracy/Specificity/
# Sensitivity) are cleared out.
$ws.Range("A8").Value = "ResNet50"
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()

# Row 9: the old placeholder cells A9:I9 are no longer needed and are
# removed entirely.
$ws.Range("A9:I9").Clear()

# Move/update the active selection.
$ws.Range("I11").Select() | Out-Null
